# Fruta / hortaliza, semanal
# Insert a new weekly price-report row (row 4) into the Cebollín sheet,
# pushing the existing rows 4-34 down to rows 5-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:34 down by one to make room for the new record.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly observation.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44847
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112037
$ws.Range("G4").Value = "Cebollín"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 7500
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7750
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 2583
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
